# The workbook's data table (Sheet1, rows 2:103) records one price entry
# per row for "Agrícola del Norte S.A. de Arica" / "Locoto". A new weekly
# entry (two quality rows: Primera + Segunda) is inserted at the top of the
# table (rows 18:19), pushing the existing rows 18:103 down to 20:105 and
# growing the used range from A1:R103 to A1:R105.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows; Excel shifts everything from row 18 down to 20,
# carrying over formatting (e.g. column D's date style) automatically.
$ws.Rows("18:19").Insert()

$newDate = Get-Date -Year 2022 -Month 5 -Day 20 -Hour 0 -Minute 0 -Second 0

# Row 18: Calidad "Primera"
$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(18, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(18, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(18, 4).Value = $newDate
$ws.Cells.Item(18, 5).Value = 15
$ws.Cells.Item(18, 6).Value = 100112042
$ws.Cells.Item(18, 7).Value = "Locoto"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 100
$ws.Cells.Item(18, 11).Value = 54000
$ws.Cells.Item(18, 12).Value = 55000
$ws.Cells.Item(18, 13).Value = 54500
$ws.Cells.Item(18, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(18, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(18, 16).Value = 2725
$ws.Cells.Item(18, 17).Value = 20
$ws.Cells.Item(18, 18).Value = "Hortaliza"

# Row 19: Calidad "Segunda"
$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(19, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(19, 4).Value = $newDate
$ws.Cells.Item(19, 5).Value = 15
$ws.Cells.Item(19, 6).Value = 100112042
$ws.Cells.Item(19, 7).Value = "Locoto"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Segunda"
$ws.Cells.Item(19, 10).Value = 160
$ws.Cells.Item(19, 11).Value = 49000
$ws.Cells.Item(19, 12).Value = 50000
$ws.Cells.Item(19, 13).Value = 49500
$ws.Cells.Item(19, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(19, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(19, 16).Value = 2475
$ws.Cells.Item(19, 17).Value = 20
$ws.Cells.Item(19, 18).Value = "Hortaliza"
